# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (a duplicate of the "2022-Q3" report
# layout, populated with the new quarter's fund-holding figures) right
# after the "总计" summary sheet, and records the new quarter in the
# "总计" summary table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet
#    (keeps header/border/font styling identical to the other quarterly
#    report tabs) and drop it right after "总计".
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Copy($null, $wb.Worksheets.Item("总计"))
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Fund-holding rows for 2022-Q4.
# code, name, size, stockPosition, positionRatio, marketValue, rank
$q4rows = @(
    @("159855", "银华中证影视主题ETF", "1.01", "97.80", "4.05", "0.0409", 9),
    @("516620", "国泰中证影视主题ETF", "0.71", "98.01", "4.29", "0.0305", 6),
    @("008778", "嘉实中证500指数增强A", "0.60", "93.52", "1.75", "0.0105", 8),
    @("008779", "嘉实中证500指数增强C", "0.40", "93.52", "1.75", "0.0070", 8)
)

for ($i = 0; $i -lt $q4rows.Length; $i++) {
    $r = $i + 2
    $data = $q4rows[$i]

    # Force text storage for the numeric-looking columns (fund code +
    # the four decimal figures) so leading/trailing zeros survive, same
    # as the other quarterly sheets.
    $q4.Range("B" + $r).NumberFormat = "@"
    $q4.Range("D" + $r + ":G" + $r).NumberFormat = "@"

    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $data[0]
    $q4.Cells.Item($r, 3).Value = $data[1]
    $q4.Cells.Item($r, 4).Value = $data[2]
    $q4.Cells.Item($r, 5).Value = $data[3]
    $q4.Cells.Item($r, 6).Value = $data[4]
    $q4.Cells.Item($r, 7).Value = $data[5]
    $q4.Cells.Item($r, 8).Value = $data[6]
}

# Rows 2-3 already carried column-A styling from the copied template;
# replicate it onto the two freshly appended rows (4-5) as well.
$q4.Cells.Item(3, 1).Copy()
$q4.Cells.Item(4, 1).PasteSpecial(-4122)
$q4.Cells.Item(5, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the 2022-Q4 totals as the
#    new first data row, pushing the existing quarters down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$totalRows = @(
    @("2022-Q4", 4, 0.09),
    @("2022-Q3", 2, 0.13),
    @("2021-Q4", 2, 0.12)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $data = $totalRows[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $data[0]
    $total.Cells.Item($r, 3).Value = $data[1]
    $total.Cells.Item($r, 4).Value = $data[2]
}
# New row 4 needs the same column-A styling as the existing data rows.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(4, 1).PasteSpecial(-4122)
